$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A13").Value = "office/stringsearch/runme_large.sh"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("A13").Select()
